# Update "想去人数" (want-to-go count) figures in column F for the two
# sheets that carry this data table ("展览" and "全部类型"). Both sheets
# contain identical rows, so the same cell->new-value map is applied to
# each.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F4"  = 119
    "F5"  = 106
    "F6"  = 475
    "F7"  = 55
    "F9"  = 599
    "F10" = 38
    "F11" = 323
    "F15" = 106
    "F16" = 17
    "F20" = 1013
    "F21" = 1422
    "F23" = 342
    "F24" = 185
    "F31" = 283
    "F32" = 1646
    "F36" = 592
    "F38" = 3790
    "F39" = 1
    "F42" = 942
    "F43" = 44
    "F46" = 80
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
